# Generate Report for Handback
# Applies the "handback" status update to the localization-status workbook:
#  - Status cells move from "Ready for handoff" -> "Handed back: in sync with en-US"
#  - The zh-cn / de-de detail sheets get their "Latest Target File" (I),
#    "Latest Handback File" (J) and "Latest Handback DateTime" (K) columns
#    filled in now that a handback has happened, including a new hyperlink
#    on the "Latest Target File" cell (mirrors the existing source-file
#    hyperlink in column A).
#  - A couple of columns are widened to comfortably fit the long file names.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status text: every cell currently showing "Ready for handoff" flips to
#    "Handed back: in sync with en-US" (Overview summary + both language
#    sheets' Status column).
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

$wsZhCn.Range("C2").Value = $handedBack
$wsZhCn.Range("C3").Value = $handedBack

$wsDeDe.Range("C2").Value = $handedBack
$wsDeDe.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I), Latest Handback File (J)
#    and Latest Handback DateTime (K) for both data rows, and add a
#    hyperlink on the target-file cell (same URL as the column-A source
#    hyperlink).
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md", [Type]::Missing, [Type]::Missing, "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md", [Type]::Missing, [Type]::Missing, "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md", [Type]::Missing, [Type]::Missing, "42647d57-8228-4722-a6e3-4fd76a0d03a6.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md", [Type]::Missing, [Type]::Missing, "42647d57-8228-4722-a6e3-4fd76a0d03a6.md") | Out-Null

$wsZhCn.Range("J2").Value = "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.9c652668270faba77ec5a33cb84f14b0ab7c5182.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "42647d57-8228-4722-a6e3-4fd76a0d03a6.8f99d4fb6fa6ab4fea68a2aee59820cebba79779.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-24 17:04:30"
$wsZhCn.Range("K3").Value = "2016-08-24 17:04:30"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same idea, with its own (slightly later) handback
#    timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md", [Type]::Missing, [Type]::Missing, "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md", [Type]::Missing, [Type]::Missing, "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md", [Type]::Missing, [Type]::Missing, "42647d57-8228-4722-a6e3-4fd76a0d03a6.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dfbf81f0c555b338ba000af619a74d12e67ac1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md", [Type]::Missing, [Type]::Missing, "42647d57-8228-4722-a6e3-4fd76a0d03a6.md") | Out-Null

$wsDeDe.Range("J2").Value = "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.9c652668270faba77ec5a33cb84f14b0ab7c5182.de-de.xlf"
$wsDeDe.Range("J3").Value = "42647d57-8228-4722-a6e3-4fd76a0d03a6.8f99d4fb6fa6ab4fea68a2aee59820cebba79779.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-24 17:04:40"
$wsDeDe.Range("K3").Value = "2016-08-24 17:04:40"

# ---------------------------------------------------------------------------
# 4. Widen the columns that now hold the long handoff/handback file names so
#    the new hyperlinked values are readable.
# ---------------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17

$wsZhCn.Range("C1").ColumnWidth = 29.17
$wsZhCn.Range("I1").ColumnWidth = 39.17
$wsZhCn.Range("J1").ColumnWidth = 39.17

$wsDeDe.Range("C1").ColumnWidth = 29.17
$wsDeDe.Range("I1").ColumnWidth = 39.17
$wsDeDe.Range("J1").ColumnWidth = 39.17
